$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

$ws.Range("A5").Select()
